$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E (shifts existing E:H to F:I), mirroring how the
# planet "Color" column was added to the cheat-sheet ahead of the ratio data.
$ws.Columns("E:E").Insert()

# The newly inserted column should carry no inherited formatting.
$ws.Range("E1:E10").ClearFormats()

$ws.Range("E1").Value = "Color"
$ws.Range("E2").Value = "Brown"
$ws.Range("E3").Value = "Pink"
$ws.Range("E4").Value = "Red"
$ws.Range("E5").Value = "Gray"
$ws.Range("E6").Value = "Green"
$ws.Range("E7").Value = "Yellow"
$ws.Range("E8").Value = "Orange"
$ws.Range("E9").Value = "Blue"
$ws.Range("E10").Value = "Purple"

$ws.Range("E10").Select()
